$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 ("erdi") ---
$ws.Range("B3").Value = 350.0
$ws.Range("C3").Value = 5.0
$ws.Range("D3").Value = 2.55
$ws.Range("E3").Value = 0.15
$ws.Range("F3").Value = 2.46
$ws.Range("G3").Value = 7.5
$ws.Range("J3").Value = 12.0

# Newly populated cells on row 3
$ws.Range("M3").Value = 0.8999999999999999
$ws.Range("N3").Value = 3.7525
$ws.Range("O3").Value = 97.625
$ws.Range("P3").Value = 90.1
$ws.Range("Q3").Value = 67.575
$ws.Range("R3").Value = 45.0

# --- Add new row 4 ("ege") ---
$ws.Range("A4").Value = "ege"
$ws.Range("B4").Value = 162.5
$ws.Range("C4").Value = 35.0
$ws.Range("D4").Value = 2.85
$ws.Range("E4").Value = 4.05
$ws.Range("F4").Value = 2.91
$ws.Range("G4").Value = 10.0
$ws.Range("H4").Value = 10.0
$ws.Range("I4").Value = 130.0
$ws.Range("J4").Value = 6.0
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.9750000000000001
$ws.Range("M4").Value = 0.44999999999999996
$ws.Range("N4").Value = 2.86375
$ws.Range("O4").Value = 23.1875
$ws.Range("P4").Value = 34.55
$ws.Range("Q4").Value = 25.912499999999994
$ws.Range("R4").Value = 15.0
